# Append two new daily-data rows (dates 45988 and 45989) to each of the six
# worksheets, matching the pattern already present in the sheet (column A =
# date serial formatted as "YYYY-MM-DD HH:MM:SS", column B = plain number;
# last row of each pair is a zero placeholder for the not-yet-reported day).

$wb = $excel.ActiveWorkbook

$dateFormat = "YYYY-MM-DD HH:MM:SS"

# sheet index (1-based) -> [ lastExistingRow, valueForDay45988 ]
$sheetInfo = @(
    @{ Index = 1; LastRow = 118; Value1 = 450515 },
    @{ Index = 2; LastRow = 118; Value1 = 59252 },
    @{ Index = 3; LastRow = 118; Value1 = 13921 },
    @{ Index = 4; LastRow = 109; Value1 = 18417 },
    @{ Index = 5; LastRow = 108; Value1 = 3727 },
    @{ Index = 6; LastRow = 108; Value1 = 1344 }
)

foreach ($info in $sheetInfo) {
    $ws = $wb.Worksheets.Item($info.Index)

    $row1 = $info.LastRow + 1
    $row2 = $info.LastRow + 2

    $ws.Cells.Item($row1, 1).Value = 45988
    $ws.Cells.Item($row1, 1).NumberFormat = $dateFormat
    $ws.Cells.Item($row1, 2).Value = $info.Value1

    $ws.Cells.Item($row2, 1).Value = 45989
    $ws.Cells.Item($row2, 1).NumberFormat = $dateFormat
    $ws.Cells.Item($row2, 2).Value = 0
}
